# Actualización desde MV -datos-
# Appends new daily rows (06-08-2021 .. 06-09-2021) to Sheet1, rows 151-172.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, date(col A), Cupo(B), Monto demandado(C), Monto total adjudicado(D),
#             Monto adjudicado bancos(E), Monto adjudicado AFP(F), Tasa(G)
# $null means the cell must be left blank (not written / not present).
$rows = @(
    @(151, "06-08-2021", 800000,  494000,  400000,  261000,  139000,  0.78),
    @(152, "09-08-2021", 900000,  1227000, 900000,  518000,  382000,  0.75),
    @(153, "10-08-2021", 1000000, 2156000, 1200000, 694000,  506000,  0.75),
    @(154, "11-08-2021", 900000,  1807000, 1260000, 913000,  347000,  0.75),
    @(155, "12-08-2021", 1200000, 2055000, 1560000, 1060000, 500000,  0.75),
    @(156, "13-08-2021", 900000,  2540000, 1350000, 500000,  850000,  0.73),
    @(157, "16-08-2021", 1200000, 2961000, 1800000, 1280000, 520000,  0.73),
    @(158, "17-08-2021", 1200000, 2772000, 1800000, 753000,  1047000, 0.73),
    @(159, "18-08-2021", 1100000, 2512000, 1650000, 833000,  817000,  0.75),
    @(160, "19-08-2021", 1100000, 2578000, 1650000, 1057000, 593000,  0.75),
    @(161, "20-08-2021", 700000,  1936000, 1050000, 784000,  266000,  0.75),
    @(162, "23-08-2021", 1500000, 2054000, 1500000, 1246000, 254000,  0.88),
    @(163, "24-08-2021", 1500000, 3018000, 2250000, 1462000, 788000,  0.89),
    @(164, "25-08-2021", 1500000, 1694000, 750000,  516000,  234000,  0.89),
    @(165, "26-08-2021", 1600000, 1945000, 1600000, 1095000, 505000,  0.99),
    @(166, "27-08-2021", 1600000, 1506000, 1506000, 1010000, 496000,  1.05),
    @(167, "30-08-2021", 1600000, $null,   0,       $null,   $null,   $null),
    @(168, "31-08-2021", 2000000, $null,   0,       $null,   $null,   $null),
    @(169, "01-09-2021", 2000000, 3222000, 3000000, 2443000, 557000,  1.49),
    @(170, "02-09-2021", 2000000, 3579000, 3000000, 1596000, 1404000, 1.49),
    @(171, "03-09-2021", 2000000, 2719000, 2000000, 1446000, 554000,  1.47),
    @(172, "06-09-2021", 2200000, 2821000, 2821000, 1850000, 971000,  1.5)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $dateText = $r[1]

    $cellA = $ws.Cells.Item($rowNum, 1)
    # Some dd-mm-yyyy strings (day <= 12) are ambiguous and would otherwise be
    # auto-parsed into date serials by Excel's normal text-entry heuristics
    # (e.g. "06-08-2021" -> 08-Jun-2021). Force those to remain plain text,
    # matching how the rest of the date column is stored.
    if ($dateText -match '^(0[1-9]|1[0-2])-\d{2}-\d{4}$') {
        $cellA.NumberFormat = "@"
    }
    $cellA.Value = $dateText

    for ($col = 2; $col -le 7; $col++) {
        $val = $r[$col]
        if ($null -ne $val) {
            $ws.Cells.Item($rowNum, $col).Value = $val
        }
    }
}
